$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.048278093338013
$ws.Range("B1").Value = 3.456063985824585
$ws.Range("C1").Value = 3.414569139480591
$ws.Range("D1").Value = 2.017139434814453
$ws.Range("E1").Value = 1.162735819816589
